$d = $word.ActiveDocument

# The outline originally repeated its whole bullet list a second time after the
# "Conclusion & Discussion" line. That whole duplicated block (from right after
# "Conclusion & Discussion" through the very end of the document) is removed,
# leaving "Conclusion & Discussion" as the final paragraph.

$findRange = $d.Content
$found = $findRange.Find.Execute("Conclusion & Discussion", $false, $false, $false, `
                                  $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Grow the found range to cover the whole paragraph, including its paragraph mark.
    $findRange.Expand(4) | Out-Null   # wdParagraph

    $paraEnd = $findRange.End
    $bodyEnd = $d.Content.End

    if ($bodyEnd -gt $paraEnd) {
        $delRange = $d.Range($paraEnd, $bodyEnd)
        $delRange.Delete()
    }
}
